$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Swap the Owner Name values between row 3 and row 4
$e3 = $ws.Range("E3").Value2
$e4 = $ws.Range("E4").Value2
$ws.Range("E3").Value2 = $e4
$ws.Range("E4").Value2 = $e3
